$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the enemy task: "Work on orbs" -> "Work On Menu"
$ws.Range("B7").Value = "Work On Menu"

# Move the active selection to B7 to match the saved view state
$ws.Range("B7").Select()
